$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Update "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:47 AM"

# Update Foreman name
$ws.Range("G8").Value = "Jonathan Estrada"

# Update Billing Period
$ws.Range("C10").Value = "06/30/2025 to 07/06/25"

# Update Job #
$ws.Range("G13").Value = "709-1"
